$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getDataGraphQL")

# --- Fill column A for rows 12-23 (in row order) ---
$ws.Range("A12").Value = 'iems-config-CIMSOURCE_TRAINCONFIG'
$ws.Range("A13").Value = 'iems-config-CIMSOURCE_TRAINRESULT'
$ws.Range("A14").Value = 'iems-config-CIMSOURCE_PREDICTCONFIG'
$ws.Range("A15").Value = 'iems-config-CIMSOURCE_PREDICTRESULT'
$ws.Range("A16").Value = 'iems-config-CIMSOURCE_OPTIMALIZERESULT'
$ws.Range("A17").Value = 'iems-config-CIMSOURCE_OPTIMALIZECONFIG'
$ws.Range("A18").Value = 'iems-config-CIMSOURCE_SIMULATION'
$ws.Range("A19").Value = 'iems-config-CIMSOURCE_TOPOLOGY'
$ws.Range("A20").Value = 'iems-config-CIMSOURCE_OPTIMIZE_TARGET_TYPE'
$ws.Range("A21").Value = 'iems-config-CIMSOURCE_DASHBOARDCONFIG'
$ws.Range("A22").Value = 'iems-config-CIMSOURCE_METADATA'
$ws.Range("A23").Value = 'iems-config-CIMSOURCE_SECTIONALIZEDLINEARIZATION'

# --- Fill column B, D, E, F for rows 12-23 (in row order) ---
$ws.Range("B12").Value = "good request, data retrieved"
$ws.Range("D12").Value = 200
$ws.Range("E12").Value = 100000
$ws.Range("F12").Value = "Successfully"
$ws.Range("B13").Value = "good request, data retrieved"
$ws.Range("D13").Value = 200
$ws.Range("E13").Value = 100000
$ws.Range("F13").Value = "Successfully"
$ws.Range("B14").Value = "good request, data retrieved"
$ws.Range("D14").Value = 200
$ws.Range("E14").Value = 100000
$ws.Range("F14").Value = "Successfully"
$ws.Range("B15").Value = "good request, data retrieved"
$ws.Range("D15").Value = 200
$ws.Range("E15").Value = 100000
$ws.Range("F15").Value = "Successfully"
$ws.Range("B16").Value = "good request, data retrieved"
$ws.Range("D16").Value = 200
$ws.Range("E16").Value = 100000
$ws.Range("F16").Value = "Successfully"
$ws.Range("B17").Value = "good request, data retrieved"
$ws.Range("D17").Value = 200
$ws.Range("E17").Value = 100000
$ws.Range("F17").Value = "Successfully"
$ws.Range("B18").Value = "good request, data retrieved"
$ws.Range("D18").Value = 200
$ws.Range("E18").Value = 100000
$ws.Range("F18").Value = "Successfully"
$ws.Range("B19").Value = "good request, data retrieved"
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 100000
$ws.Range("F19").Value = "Successfully"
$ws.Range("B20").Value = "good request, data retrieved"
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = 100000
$ws.Range("F20").Value = "Successfully"
$ws.Range("B21").Value = "good request, data retrieved"
$ws.Range("D21").Value = 200
$ws.Range("E21").Value = 100000
$ws.Range("F21").Value = "Successfully"
$ws.Range("B22").Value = "good request, data retrieved"
$ws.Range("D22").Value = 200
$ws.Range("E22").Value = 100000
$ws.Range("F22").Value = "Successfully"
$ws.Range("B23").Value = "good request, data retrieved"
$ws.Range("D23").Value = 200
$ws.Range("E23").Value = 100000
$ws.Range("F23").Value = "Successfully"

# --- Fill column C for rows 12,13,15,16,17,18,19,20,21,22,23 (skip 14) then 14 last ---
$ws.Range("C12").Value = '{CIMSOURCE_TRAINCONFIG(cond:"",order:"") { createTime optimize_target_type_id simu_id train_cfg_ANN_YN train_cfg_LR_YN train_cfg_RFR_YN train_cfg_SVR_YN train_cfg_SmpEndTime train_cfg_SmpStTime train_cfg_Tree_YN train_cfg_aDa_YN train_cfg_obj_PT train_cfg_obligate1 train_cfg_obligate2 train_cfg_obligate3 train_cfg_obligate4 train_cfg_obligate5 train_cfg_obligate6 train_cfg_predict train_cfg_sample_PT train_cfg_time train_cfg_timetrg_YN}}'
$ws.Range("C13").Value = '{CIMSOURCE_TRAINRESULT(cond:"",order:"") { variables train_result_obligate4 train_result_obligate5 train_result_obligate2 train_result_obligate3 optimize_target_type_id simu_id train_result_obligate1 config_id train_result_obligate6 runid }}'
$ws.Range("C15").Value = '{CIMSOURCE_PREDICTRESULT(cond:"",order:"config_id ASC") { Predict_result_time Predict_result_aDa Predict_result_ANN optimize_target_type_id Predict_result_SP5 Predict_result_Tree Predict_result_LR Predict_result_SVR runid createTime variables predict_train_obj_wgt_defval Predict_result_obligate2 Predict_result_obligate1 Predict_result_RFR Predict_result_obligate6 simu_id Predict_result_obligate5 Predict_result_obligate4 Predict_result_obligate3 config_id Predict_result_AVG } }'
$ws.Range("C16").Value = '{CIMSOURCE_PREDICTRESULT(cond:"",order:"config_id ASC") { Predict_result_time Predict_result_aDa Predict_result_ANN optimize_target_type_id Predict_result_SP5 Predict_result_Tree Predict_result_LR Predict_result_SVR runid createTime variables predict_train_obj_wgt_defval Predict_result_obligate2 Predict_result_obligate1 Predict_result_RFR Predict_result_obligate6 simu_id Predict_result_obligate5 Predict_result_obligate4 Predict_result_obligate3 config_id Predict_result_AVG } }'
$ws.Range("C17").Value = '{CIMSOURCE_OPTIMALIZECONFIG(cond:"",order:"") { isActive optimize_target_type_id optimalizeType opt_cfg_opt_starttime opt_cfg_user_def opt_pred_model opt_cfg_obligate3 rollTasksId opt_cfg_obligate2 opt_cfg_obligate5 createTime opt_cfg_obligate4 opt_cfg_obligate6 opt_cfg_emi_frac opt_name dash_cfg_id simu_id opt_cfg_opt_day opt_cfg_opt_frequency isRTCOP opt_cfg_cost_frac opt_methods rollType opt_cfg_obligate1 } }'
$ws.Range("C18").Value = '{CIMSOURCE_SIMULATION(cond:"",order:"") { updatatime historymap chartsetting createtime rtcountmap hiscountmap isTraining obligate1 obligate2 channels defcountmap user_id ispredict name obligate5 obligate6 treefile obligate3 obligate4 } }'
$ws.Range("C19").Value = '{CIMSOURCE_TOPOLOGY(cond:"",order:"") { create_time topo_name simu_id obligate1 obligate2 topo_xml update_time update_user parent_id obligate5 obligate6 obligate3 obligate4 } }'
$ws.Range("C20").Value = '{CIMSOURCE_OPTIMIZE_TARGET_TYPE(cond:"",order:"") { obligate1 obligate2 optimize_topo_selection optimize_type optimize_name obligate5 optimize_data_type obligate3 obligate4 } }'
$ws.Range("C21").Value = '{CIMSOURCE_DASHBOARDCONFIG(cond:"",order:"") { dashparam name simu_id optimize_target_type_id } }'
$ws.Range("C22").Value = '{CIMSOURCE_METADATA(cond:"",order:"") { obligate1 obligate2 name obligate5 params category obligate3 obligate4 } }'
$ws.Range("C23").Value = '{CIMSOURCE_SECTIONALIZEDLINEARIZATION(cond:"",order:"") { tensor_index_counts target obligate1 obligate2 device_name model_name param obligate5 obligate3 obligate4 tensor_data } }'
$ws.Range("C14").Value = '{CIMSOURCE_PREDICTCONFIG(cond:"",order:"") { predict_train_model createTime predict_train_obj_wgt_defval predict_cfg_obligate6 predict_cfg_obligate5 predict_cfg_obligate4 predict_cfg_obligate3 predict_cfg_obligate2 predict_cfg_obligate1 optimize_target_type_id simu_id predict_cfg_period predict_cfg_day predict_cfg_time } }'

# --- Update selection ---
$ws.Range("B27").Select() | Out-Null
